$d = $word.ActiveDocument

# 1. Update the delivery date text from 14/12 to 13/09
$d.Content.Find.Execute("14/12", $true, $false, $false, $false, $false,
                         $true, 1, $false, "13/09", 2) | Out-Null

# 2. Remove the leftover "_GoBack" bookmark (last-edit marker) near "Listas, Tuplas"
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
